$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 11629263
$ws.Range("I15").Value = 11629263
$ws.Range("K15").Value = 34887789
$ws.Range("M15").Value = -34887620
$ws.Range("H28").Value = 909.1667
$ws.Range("I28").Value = 407.7143
$ws.Range("K28").Value = 407.7143
$ws.Range("M28").Value = 77.28570000000002
$ws.Range("H74").Value = 5479.6
$ws.Range("J74").Value = 6166.3335
$ws.Range("L74").Value = 6166.3335
$ws.Range("N74").Value = -8038.3335
$ws.Range("H77").Value = 5479.6
$ws.Range("J77").Value = 6166.3335
$ws.Range("L77").Value = 30831.6675
$ws.Range("N77").Value = -40191.6675
$ws.Range("H107").Value = 480.42856
$ws.Range("J107").Value = 553
$ws.Range("L107").Value = 553
$ws.Range("N107").Value = -4393
$ws.Range("H111").Value = 2700.6667
$ws.Range("I111").Value = 2673.5
$ws.Range("J111").Value = 2722.4
$ws.Range("K111").Value = 8020.5
$ws.Range("L111").Value = 8167.200000000001
$ws.Range("M111").Value = -4953.5
$ws.Range("N111").Value = -14301.2
$ws.Range("H113").Value = 4332.3335
$ws.Range("I113").Value = 3997
$ws.Range("K113").Value = 3997
$ws.Range("M113").Value = -743
$ws.Range("H116").Value = 34641884
$ws.Range("I116").Value = 50702908
$ws.Range("J116").Value = 25005272
$ws.Range("K116").Value = 50702908
$ws.Range("L116").Value = 25005272
$ws.Range("M116").Value = -50699466
$ws.Range("N116").Value = -25012156

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 27401.666
$ws.Range("I61").Value = 31982.2
$ws.Range("K61").Value = 31982.2
$ws.Range("M61").Value = -31770.2
$ws.Range("H74").Value = 13890850
$ws.Range("I74").Value = 27779222
$ws.Range("K74").Value = 27779222
$ws.Range("M74").Value = -27778348
$ws.Range("H77").Value = 13890850
$ws.Range("I77").Value = 27779222
$ws.Range("K77").Value = 138896110
$ws.Range("M77").Value = -138891742
$ws.Range("H110").Value = 930494.8
$ws.Range("I110").Value = 1362325.5
$ws.Range("K110").Value = 1362325.5
$ws.Range("M110").Value = -1360280.5
$ws.Range("H136").Value = 27401.666
$ws.Range("I136").Value = 31982.2
$ws.Range("K136").Value = 95946.60000000001
$ws.Range("M136").Value = -93396.60000000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 45541.25
$ws.Range("J81").Value = 45541.25
$ws.Range("L81").Value = 45541.25
$ws.Range("N81").Value = -47663.25
$ws.Range("H84").Value = 45541.25
$ws.Range("J84").Value = 45541.25
$ws.Range("L84").Value = 136623.75
$ws.Range("N84").Value = -147231.75
$ws.Range("H86").Value = 2278
$ws.Range("J86").Value = 2427.3333
$ws.Range("L86").Value = 2427.3333
$ws.Range("N86").Value = -4673.3333
$ws.Range("H89").Value = 2278
$ws.Range("J89").Value = 2427.3333
$ws.Range("L89").Value = 12136.6665
$ws.Range("N89").Value = -23368.6665
$ws.Range("H107").Value = 1725.5714
$ws.Range("I107").Value = 1729.1666
$ws.Range("J107").Value = 1704
$ws.Range("K107").Value = 1729.1666
$ws.Range("L107").Value = 1704
$ws.Range("M107").Value = 190.8334
$ws.Range("N107").Value = -5544
$ws.Range("H134").Value = 3462.4546
$ws.Range("I134").Value = 1555.5
$ws.Range("K134").Value = 4666.5
$ws.Range("M134").Value = -2131.5
$ws.Range("H138").Value = 84140
$ws.Range("J138").Value = 84140
$ws.Range("L138").Value = 84140
$ws.Range("N138").Value = -94420

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1634.5714
$ws.Range("I16").Value = 1657
$ws.Range("K16").Value = 1657
$ws.Range("M16").Value = -1370
$ws.Range("H58").Value = 1432571
$ws.Range("I58").Value = 1670499.5
$ws.Range("J58").Value = 5000
$ws.Range("K58").Value = 1670499.5
$ws.Range("L58").Value = 5000
$ws.Range("M58").Value = -1670296.5
$ws.Range("N58").Value = -5406
$ws.Range("H105").Value = 3789210.5
$ws.Range("J105").Value = 2000
$ws.Range("L105").Value = 2000
$ws.Range("N105").Value = -5494
$ws.Range("H113").Value = 1634.5714
$ws.Range("I113").Value = 1657
$ws.Range("K113").Value = 1657
$ws.Range("M113").Value = 513
$ws.Range("H136").Value = 1432571
$ws.Range("I136").Value = 1670499.5
$ws.Range("K136").Value = 5011498.5
$ws.Range("M136").Value = -5008948.5
$ws.Range("N136").Value = -20100

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 286.2
$ws.Range("I86").Value = 316.66666
$ws.Range("J86").Value = 240.5
$ws.Range("K86").Value = 949.9999799999999
$ws.Range("L86").Value = 721.5
$ws.Range("M86").Value = 236.0000200000001
$ws.Range("N86").Value = -3093.5
$ws.Range("H89").Value = 286.2
$ws.Range("I89").Value = 316.66666
$ws.Range("J89").Value = 240.5
$ws.Range("K89").Value = 2849.99994
$ws.Range("L89").Value = 2164.5
$ws.Range("M89").Value = 3078.00006
$ws.Range("N89").Value = -14020.5
$ws.Range("H129").Value = 1951.3334
$ws.Range("J129").Value = 1677
$ws.Range("L129").Value = 5031
$ws.Range("N129").Value = -15031

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 15943968
$ws.Range("I102").Value = 22180246
$ws.Range("K102").Value = 22180246
$ws.Range("M102").Value = -22178624
$ws.Range("H107").Value = 2976808.2
$ws.Range("I107").Value = 11905404
$ws.Range("J107").Value = 609.75
$ws.Range("K107").Value = 11905404
$ws.Range("L107").Value = 609.75
$ws.Range("M107").Value = -11903484
$ws.Range("N107").Value = -4449.75
$ws.Range("H113").Value = 2433
$ws.Range("I113").Value = 2549.5
$ws.Range("J113").Value = 2200
$ws.Range("K113").Value = 2549.5
$ws.Range("L113").Value = 2200
$ws.Range("M113").Value = -379.5
$ws.Range("N113").Value = -6540

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4499.25
$ws.Range("I61").Value = 4499.5
$ws.Range("J61").Value = 4499
$ws.Range("K61").Value = 4499.5
$ws.Range("L61").Value = 4499
$ws.Range("M61").Value = -4297.5
$ws.Range("N61").Value = -4903
$ws.Range("H113").Value = 4499.25
$ws.Range("I113").Value = 4499.5
$ws.Range("J113").Value = 4499
$ws.Range("K113").Value = 4499.5
$ws.Range("L113").Value = 4499
$ws.Range("M113").Value = -2329.5
$ws.Range("N113").Value = -8839

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H22").Value = 14995
$ws.Range("J22").Value = 14995
$ws.Range("L22").Value = 14995
$ws.Range("N22").Value = -15581
$ws.Range("H100").Value = 910614
$ws.Range("I100").Value = 1334242.8
$ws.Range("K100").Value = 2668485.6
$ws.Range("M100").Value = -2667944.6
$ws.Range("H113").Value = 908.4167
$ws.Range("J113").Value = 1136.6666
$ws.Range("L113").Value = 3409.9998
$ws.Range("N113").Value = -7749.9998
$ws.Range("H136").Value = 8979.361999999999
$ws.Range("I136").Value = 3650.963
$ws.Range("J136").Value = 12404.762
$ws.Range("K136").Value = 10952.889
$ws.Range("L136").Value = 37214.286
$ws.Range("M136").Value = -8402.889000000001
